# Update with new legend test
#
# Reproduces:
#   - workbook window geometry change (best-effort; Excel's own remembered
#     OS window rect has no effect on workbook contents)
#   - a new "VerifyLegend" worksheet after "VerifyMapLayers" with its legend
#     icon/alt/src reference table
#   - the "Login Info" sheet's remembered selection moving from B11 to D3

$wb = $excel.ActiveWorkbook

# --- cosmetic: remembered workbook window position/size -------------------
[void]($excel.ActiveWindow.Left   = 12060)
[void]($excel.ActiveWindow.Top    = 12400)
[void]($excel.ActiveWindow.Width  = 18840)
[void]($excel.ActiveWindow.Height = 10140)

# --- add the new "VerifyLegend" sheet, placed after "VerifyMapLayers" -----
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$legend = $wb.Worksheets.Add($null, $lastSheet)
$legend.Name = "VerifyLegend"

# Fill cells in the same order the original author must have used so that
# new shared-string entries land at the same indices the recorded diff
# expects (icon names down column A first, then the header row, then the
# "icon" alt-text column, before the rows that just repeat existing text).
$legend.Range("A2").Value = "Closure"
$legend.Range("A3").Value = "Restriction"
$legend.Range("A4").Value = "Crash"
$legend.Range("A5").Value = "Warning"
$legend.Range("A6").Value = "Info"
$legend.Range("A7").Value = "Coming Soon"

$legend.Range("A1").Value = "Legend Icons:"
$legend.Range("B1").Value = "Alts:"
$legend.Range("C1").Value = "Img srcs:"

$legend.Range("B2").Value = "Closure icon"
$legend.Range("B3").Value = "Restriction icon"

$legend.Range("B4").Value = "Crash"
$legend.Range("B5").Value = "Warning"
$legend.Range("B6").Value = "Info"
$legend.Range("B7").Value = "Coming Soon"

$legend.Range("A8").Value = "Construction"
$legend.Range("B8").Value = "Construction"

$legend.Range("A9").Value = "Weather Warnings"
$legend.Range("B9").Value = "Weather Warnings"

# Columns A:B best-fit to their text (closest width this host can reproduce).
$legend.Columns("A:B").ColumnWidth = 15.428571428571429

# Leave the new sheet's remembered selection on B4.
[void]$legend.Range("B4").Select()

# --- "Login Info" sheet: move remembered selection from B11 to D3 ---------
$loginInfo = $wb.Worksheets.Item("Login Info")
[void]$loginInfo.Range("D3").Select()

# Restore "Login Info" as the active/selected tab (it was tabSelected before
# the edit and stays that way afterwards).
[void]$loginInfo.Activate()
